$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 2156.039195262611
$ws.Range("D3").Value = 158.4280988661488

$ws.Range("B4").Value = 986.4419444683325
$ws.Range("D4").Value = 220.1866249282304

$ws.Range("B5").Value = 45.003

$ws.Range("B6").Value = 1420

$ws.Range("B7").Value = 2015.0275
$ws.Range("D7").Value = 160

$ws.Range("B8").Value = 2936.137249999998
$ws.Range("D8").Value = 280

$ws.Range("B9").Value = 6321.075999999999
$ws.Range("D9").Value = 1820

$ws.Range("F10").Value = 1133214201.03

$ws.Range("G11").Value = 0.7358194888566573

$ws.Range("F12").Value = 83269808.76400003
$ws.Range("G12").Value = 0.07348108476607029

$ws.Range("G13").Value = 0.1906994263772724
